$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Step 1: remove the _GoBack bookmark pair from the end of the
# "TOTAL WEEKLY TIME SPENT" row's time cell (currently after the
# " 30mins" run).
# ------------------------------------------------------------------
$totalCell  = $d.Tables.Item(2).Cell(5, 2)
$totalPara  = $totalCell.Range.Paragraphs.Item(1)
$totalRange = $totalPara.Range

$xmlNoBookmark = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w:rsidR="009951E4" w:rsidRDefault="00F77EAA"><w:pPr><w:spacing w:before="120"/><w:rPr><w:rFonts w:ascii="Ebrima" w:eastAsia="Ebrima" w:hAnsi="Ebrima" w:cs="Ebrima"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Ebrima" w:eastAsia="Ebrima" w:hAnsi="Ebrima" w:cs="Ebrima"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>1</w:t></w:r><w:r w:rsidR="00223C6A"><w:rPr><w:rFonts w:ascii="Ebrima" w:eastAsia="Ebrima" w:hAnsi="Ebrima" w:cs="Ebrima"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>hrs</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Ebrima" w:eastAsia="Ebrima" w:hAnsi="Ebrima" w:cs="Ebrima"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> 30mins</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$totalRange.InsertXML($xmlNoBookmark)

# ------------------------------------------------------------------
# Step 2: add an "N/A" run (bold, Ebrima, 10pt) to the empty
# "3.ISSUES/PROBLEMS" value cell, followed by the _GoBack bookmark
# pair (now reusing id 1, since it was freed up in step 1).
# ------------------------------------------------------------------
$issuesCell  = $d.Tables.Item(2).Cell(11, 1)
$issuesPara  = $issuesCell.Range.Paragraphs.Item(1)
$issuesRange = $issuesPara.Range

$xmlWithBookmark = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w:rsidR="009951E4" w:rsidRDefault="009951E4"><w:pPr><w:spacing w:line="720" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Ebrima" w:eastAsia="Ebrima" w:hAnsi="Ebrima" w:cs="Ebrima"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Ebrima" w:eastAsia="Ebrima" w:hAnsi="Ebrima" w:cs="Ebrima"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>N/A</w:t></w:r><w:bookmarkStart w:id="1" w:name="_GoBack"/><w:bookmarkEnd w:id="1"/></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$issuesRange.InsertXML($xmlWithBookmark)

Write-Output "Edit applied."
